# HojeSomaSE.xlsx - consolidate the second summary table (columns L:N) one
# step to the left (K:L) and drop its redundant "Lucia" sub-column, while
# also dropping the redundant "Renan" sub-column from the first summary
# table (columns H:J -> H:I). Formulas in the second table are rewritten
# from the SUMIFS-against-a-mirror-header style to a single SUMPRODUCT that
# references the (now adjacent) label column directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- 0. Stash the original header-row formatting so re-merging doesn't
#        introduce new border-split styles -------------------------------
$ws.Range("H1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L1").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # xlPasteFormats

# --- 1. Move the second table's label column (old L) to K, keeping style ---
$ws.Range("L1:L5").Copy($ws.Range("K1:K5"))

# --- 2. Move the second table's "Renan" column (old N) to L, keeping style ---
$ws.Range("N1:N5").Copy($ws.Range("L1:L5"))

# --- 3. Rewrite the L3:L5 formulas to reference the new K column/L2 header ---
$ws.Range("L3").Formula = '=SUMPRODUCT(($B$3:$B$14=K3)*($C$3:$C$14=$L$2)*($D$3:$D$14)*($E$3:$E$14))'
$ws.Range("L4").Formula = '=SUMPRODUCT(($B$3:$B$14=K4)*($C$3:$C$14=$L$2)*($D$3:$D$14)*($E$3:$E$14))'
$ws.Range("L5").Formula = '=SUMPRODUCT(($B$3:$B$14=K5)*($C$3:$C$14=$L$2)*($D$3:$D$14)*($E$3:$E$14))'

# --- 4. Drop the now-redundant columns: J (Renan, table 1) and M:N (old table 2 leftovers) ---
# NOTE: cleared cell-by-cell on purpose. Clear()/ClearContents() on a
# multi-row range whose top row crosses one of row 1's merged header cells
# only ever ends up touching the first cell of that range, so every cell is
# cleared individually here to be safe.
$ws.Range("J1").Clear()
$ws.Range("J2").Clear()
$ws.Range("J3").Clear()
$ws.Range("J4").Clear()
$ws.Range("J5").Clear()
$ws.Range("M1").Clear()
$ws.Range("M2").Clear()
$ws.Range("M3").Clear()
$ws.Range("M4").Clear()
$ws.Range("M5").Clear()
$ws.Range("N1").Clear()
$ws.Range("N2").Clear()
$ws.Range("N3").Clear()
$ws.Range("N4").Clear()
$ws.Range("N5").Clear()

# --- 5. Fix the merged header cells: H1:J1 -> H1:I1, (old) L1:N1 -> K1:L1 ---
$ws.Range("H1:J1").UnMerge()
$ws.Range("H1:I1").Merge()
$ws.Range("Z1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("L1:N1").UnMerge()
$ws.Range("K1:L1").Merge()
$ws.Range("Z2").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Z1").Clear()
$ws.Range("Z2").Clear()

# --- 6. Column widths for the now-adjacent columns ---
$ws.Columns.Item(9).ColumnWidth = 11.42578125
$ws.Columns.Item(11).ColumnWidth = 9.7109375
$ws.Columns.Item(12).ColumnWidth = 11.7109375

# --- 7. Selection / active cell ---
$ws.Range("L3").Select()

$excel.CutCopyMode = 0
$wb.Save()
